$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.554.17"
$ws.Range("E2").Value = "  +5.08%  "

$ws.Range("D3").Value = "1.721.71"
$ws.Range("E3").Value = "  +3.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5408"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2760"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06785"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.57%  "

$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.699"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.22%  "

$ws.Range("D13").Value = "1.718.44"
$ws.Range("E13").Value = "  +3.21%  "

$ws.Range("D14").Value = "1.956.91"
$ws.Range("E14").Value = "  +3.91%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5989"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.85%  "

$ws.Range("D16").Value = "0.0₅8357"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.38%  "

$ws.Range("D18").Value = "27.498.56"
$ws.Range("E18").Value = "  +4.90%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.813"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.22%  "

$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "211.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.79%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.94%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.218"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1246"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.429"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.96%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.625"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.37%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05586"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.313"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.673"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.523"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.631"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9756"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.98%  "

$ws.Range("E36").Value = "  +1.77%  "

$ws.Range("E37").Value = "  +1.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5903"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01647"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.862"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.81%  "

$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("D42").Value = "1.038.33"
$ws.Range("E42").Value = "  +0.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8360"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.58%  "

$ws.Range("D45").Value = "1.862.28"
$ws.Range("E45").Value = "  +3.83%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈116"
$ws.Range("E46").Value = "  +10.15%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "59.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.195"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4419"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.59%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.003"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05273"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.01%  "

